# Update dashboards - 2025-12-09
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - ADP Total NonFarm Private (M/M Delta)
$ws.Range("Q5").Value = -32000

# Row 29 - 5yr, 5yr Forward
$ws.Range("N29").Value = 45999
$ws.Range("Q29").Value = 2.2
$ws.Range("R29").Value = 2.18
$ws.Range("T29").Value = 2.19
$ws.Range("U29").Value = 2.18

# Row 30 - 10yr TIPS
$ws.Range("N30").Value = 45999
$ws.Range("S30").Value = 2.26

# Row 39 - Nominal Broad US Dollar Index
$ws.Range("N39").Value = 45996
$ws.Range("N39").Style = $ws.Range("N47").Style
$ws.Range("Q39").Value = 121.0615
$ws.Range("R39").Value = 121.0614
$ws.Range("S39").Value = 121.1131
$ws.Range("T39").Value = 121.5149
$ws.Range("U39").Value = 121.3615

# Row 47 - FFR
$ws.Range("N47").Value = 45996

# Row 48 - 2y UST
$ws.Range("N48").Value = 45996
$ws.Range("Q48").Value = 3.56
$ws.Range("R48").Value = 3.52
$ws.Range("S48").Value = 3.49
$ws.Range("T48").Value = 3.51
$ws.Range("U48").Value = 3.54

# Row 49 - 5y UST
$ws.Range("N49").Value = 45996
$ws.Range("Q49").Value = 3.72
$ws.Range("R49").Value = 3.68
$ws.Range("S49").Value = 3.62
$ws.Range("T49").Value = 3.66
$ws.Range("U49").Value = 3.67

# Row 50 - 10y UST
$ws.Range("N50").Value = 45996
$ws.Range("Q50").Value = 4.14
$ws.Range("R50").Value = 4.11
$ws.Range("S50").Value = 4.06
$ws.Range("T50").Value = 4.09
$ws.Range("U50").Value = 4.09

# Row 52 - BAA
$ws.Range("N52").Value = 45996
$ws.Range("Q52").Value = 5.88
$ws.Range("R52").Value = 5.87
$ws.Range("S52").Value = 5.83
$ws.Range("T52").Value = 5.85
$ws.Range("U52").Value = 5.87
